$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11, pushing existing rows 11-34 down to 12-35
$ws.Rows.Item(11).Insert()

# Populate the newly inserted row 11 with the new record
$ws.Range("A11").Value = 5
$ws.Range("B11").Value = "Macroferia Regional de Talca"
$ws.Range("C11").Value = "Maule"
$ws.Range("D11").Value = 44481
$ws.Range("E11").Value = 7
$ws.Range("F11").Value = 100112026
$ws.Range("G11").Value = "Haba"
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 500
$ws.Range("K11").Value = 8000
$ws.Range("L11").Value = 8000
$ws.Range("M11").Value = 8000
$ws.Range("N11").Value = "$/saco 25 kilos"
$ws.Range("O11").Value = "Región de O'Higgins"
$ws.Range("P11").Value = 320
$ws.Range("Q11").Value = 25
$ws.Range("R11").Value = "Hortaliza"
